$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$changes = @{
    5  = -2
    6  = 1
    14 = 4
    24 = -8
    25 = -5
    28 = 0
    30 = -4
    31 = 1
    34 = -5
    37 = -1
    41 = -1
    42 = -5
    58 = -9
    59 = -1
    61 = 0
    71 = 3
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
